$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "26.799.23"
$ws.Range("E2").Value2 = "  -1.81%  "
$ws.Range("D3").Value2 = "1.548.18"
$ws.Range("E3").Value2 = "  -1.87%  "
$ws.Range("E4").Value2 = "  +0.04%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "204.61"
$cell.Style = $origStyle
$ws.Range("E5").Value2 = "  -1.70%  "
$ws.Range("E6").Value2 = "  -1.76%  "
$ws.Range("E7").Value2 = "  +0.06%  "
$ws.Range("B8").Value2 = "Solana"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "21.41"
$cell.Style = $origStyle
$ws.Range("E8").Value2 = "  -4.04%  "
$ws.Range("B9").Value2 = "Cardano"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "0.245"
$cell.Style = $origStyle
$ws.Range("E9").Value2 = "  -1.12%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "0.0581"
$cell.Style = $origStyle
$ws.Range("E10").Value2 = "  -1.87%  "
$ws.Range("E11").Value2 = "  -1.06%  "
$ws.Range("D12").Value2 = "1.769.14"
$ws.Range("E12").Value2 = "  -1.78%  "
$ws.Range("D13").Value2 = "1.546.50"
$ws.Range("E13").Value2 = "  -2.02%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "3.68"
$cell.Style = $origStyle
$ws.Range("E14").Value2 = "  -2.78%  "
$ws.Range("E15").Value2 = "  -1.97%  "
$ws.Range("D16").Value2 = "26.786.92"
$ws.Range("E16").Value2 = "  -1.83%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "60.87"
$cell.Style = $origStyle
$ws.Range("E17").Value2 = "  -2.69%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "213.89"
$cell.Style = $origStyle
$ws.Range("E18").Value2 = "  -0.84%  "
$ws.Range("E19").Value2 = "  -1.45%  "
$ws.Range("D20").Value2 = "0.0₃0682"
$ws.Range("E20").Value2 = "  -0.86%  "
$ws.Range("E21").Value2 = "  +0.03%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "4.07"
$cell.Style = $origStyle
$ws.Range("E22").Value2 = "  -1.78%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "9.03"
$cell.Style = $origStyle
$ws.Range("E23").Value2 = "  -4.27%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "152.70"
$cell.Style = $origStyle
$ws.Range("E25").Value2 = "  +0.60%  "
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "6.51"
$cell.Style = $origStyle
$ws.Range("E26").Value2 = "  -2.71%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "14.91"
$cell.Style = $origStyle
$ws.Range("E27").Value2 = "  -0.44%  "
$ws.Range("E28").Value2 = "  +0.05%  "
$ws.Range("E29").Value2 = "  -2.18%  "
$ws.Range("E30").Value2 = "  -0.70%  "
$ws.Range("E31").Value2 = "  -3.65%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "3.17"
$cell.Style = $origStyle
$ws.Range("E32").Value2 = "  -0.55%  "
$ws.Range("D33").Value2 = "1.353.98"
$ws.Range("E33").Value2 = "  -4.09%  "
$ws.Range("E34").Value2 = "  -0.96%  "
$ws.Range("E35").Value2 = "  -3.96%  "
$ws.Range("E36").Value2 = "  -0.79%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "0.916"
$cell.Style = $origStyle
$ws.Range("E37").Value2 = "  -2.42%  "
$ws.Range("E38").Value2 = "  -2.21%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "0.523"
$cell.Style = $origStyle
$ws.Range("E39").Value2 = "  +0.50%  "
$ws.Range("E40").Value2 = "  -2.67%  "
$ws.Range("E41").Value2 = "  +0.05%  "
$ws.Range("E42").Value2 = "  -1.31%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "5.55"
$cell.Style = $origStyle
$ws.Range("E43").Value2 = "  +3.73%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "2.19"
$cell.Style = $origStyle
$ws.Range("E44").Value2 = "  +0.17%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "1.76"
$cell.Style = $origStyle
$ws.Range("E45").Value2 = "  -2.93%  "
$ws.Range("E46").Value2 = "  -1.77%  "
$ws.Range("E47").Value2 = "  -2.23%  "
$ws.Range("D48").Value2 = "1.683.30"
$ws.Range("E48").Value2 = "  -1.74%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value2 = "85.89"
$cell.Style = $origStyle
$ws.Range("E49").Value2 = "  -0.56%  "
$ws.Range("E50").Value2 = "  +2.91%  "
$ws.Range("D51").Value2 = "0.0₇0973"
$ws.Range("E51").Value2 = "  -1.69%  "
